$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-03-13"

# Update the row label for March to reflect the new "through" date
$ws.Range("A4").Value = "March (through 03-13)"

# Update March row (row 4) values for columns D through I
$ws.Range("D4").Value = 26
$ws.Range("E4").Value = 25
$ws.Range("F4").Value = 13
$ws.Range("G4").Value = 26
$ws.Range("H4").Value = 36
$ws.Range("I4").Value = 62

# Update Total row (row 5) values for columns D through I
$ws.Range("D5").Value = 157
$ws.Range("E5").Value = 162
$ws.Range("F5").Value = 92
$ws.Range("G5").Value = 167
$ws.Range("H5").Value = 378
$ws.Range("I5").Value = 362
